# QS003 test script update:
#   Add a new "check_page_accessibility" step (step + keyword rows) right
#   after the existing "Open the browser" / "Chrome" setup rows (i.e. before
#   the pre-existing row 7), and move the active selection to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 7:8 (pushes every following row down by two,
# inheriting the formatting of the row above exactly like Excel does).
$ws.Rows("7:8").Insert()

# New "step" row (7): description + expected result of the accessibility check.
$ws.Range("A7").Value = "step"
$ws.Range("B7").Value = "Check page is accesibility compliant"
$ws.Range("C7").Value = "The page should be accesibility compliant"

# New "keyword" row (8): the actual keyword that performs the check.
$ws.Range("A8").Value = "check_page_accessibility"
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A8").IndentLevel = 1

# Move the selection like it was left after the edit.
$ws.Range("C4").Select() | Out-Null

Write-Host "Inserted check_page_accessibility step/keyword rows"
